$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "P.Aur_1"
$ws.Range("A3").Value = "P.Aur_2"
$ws.Range("B3").Value = 150188
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0

$ws.Range("A4").Value = "Test056"
$ws.Range("B4").Value = 150056

$ws.Range("A5").Value = "Test196"
$ws.Range("B5").Value = 150196

$ws.Range("A6").Select()
